# Actualización automática 2025-06-02 14:03:43
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

$xlPasteFormats = -4122

# --- New column G ---

# Header G1 = "PRESUPUESTO", same formatting as F1
$ws.Cells.Item(1, 6).Copy()
$ws.Cells.Item(1, 7).PasteSpecial($xlPasteFormats)
$ws.Cells.Item(1, 7).Value = "PRESUPUESTO"

# Data cells G2:G6 = 0, same formatting as column F
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 6).Copy()
    $ws.Cells.Item($r, 7).PasteSpecial($xlPasteFormats)
    $ws.Cells.Item($r, 7).Value = 0
}

# Totals row G7 = 0, same formatting as F7
$ws.Cells.Item(7, 6).Copy()
$ws.Cells.Item(7, 7).PasteSpecial($xlPasteFormats)
$ws.Cells.Item(7, 7).Value = 0

# Column G width -> 17 characters
$ws.Columns.Item(7).ColumnWidth = 16.15
